$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Overview" (xl/worksheets/sheet1.xml) — columns A:G
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Insert a new row above the existing data row, pushing it down to row 3
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "27ed5c2c-acc1-4e8a-a34b-ed92d1d00444.md"
$ws.Range("B2").Value = "e2e\27ed5c2c-acc1-4e8a-a34b-ed92d1d00444.md"
$ws.Range("C2").Value = ".md"
$ws.Range("D2").Value = "'"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-28 04:39:28"
$ws.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Fix up hyperlinks: row-insert does not relocate the existing hyperlink,
# so clear them all and re-add both at their correct rows.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f74ea00f5892552ff51f6e0afefea33abf807e6/e2e/27ed5c2c-acc1-4e8a-a34b-ed92d1d00444.md", "", "", "e2e\27ed5c2c-acc1-4e8a-a34b-ed92d1d00444.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f74ea00f5892552ff51f6e0afefea33abf807e6/e2e/93ca8fcb-a7e6-49dd-96f9-8a8ef9039be9.md", "", "", "e2e\93ca8fcb-a7e6-49dd-96f9-8a8ef9039be9.md")

# Grow the table to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G3"))

# -----------------------------------------------------------------
# Sheet "zh-cn" (xl/worksheets/sheet2.xml) — columns A:P
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "27ed5c2c-acc1-4e8a-a34b-ed92d1d00444.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "False"
$ws.Range("G2").Value = "27ed5c2c-acc1-4e8a-a34b-ed92d1d00444.3f9bb167b8acd64d748426ebd885572042073a9a.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-28 04:39:24"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I2").Value = "'"
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = "'"
$ws.Range("J2").Style = "Normal"
$ws.Range("K2").Value = "0001-01-01 00:00:00"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L2").Value = "'"
$ws.Range("L2").Style = "Normal"
$ws.Range("M2").Value = "True"
$ws.Range("N2").Value = "'"
$ws.Range("N2").Style = "Normal"
$ws.Range("O2").Value = "False"
$ws.Range("P2").Value = "'"
$ws.Range("P2").Style = "Normal"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f74ea00f5892552ff51f6e0afefea33abf807e6/e2e/27ed5c2c-acc1-4e8a-a34b-ed92d1d00444.md", "", "", "27ed5c2c-acc1-4e8a-a34b-ed92d1d00444.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f74ea00f5892552ff51f6e0afefea33abf807e6/e2e/93ca8fcb-a7e6-49dd-96f9-8a8ef9039be9.md", "", "", "93ca8fcb-a7e6-49dd-96f9-8a8ef9039be9.md")

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))

# -----------------------------------------------------------------
# Sheet "de-de" (xl/worksheets/sheet3.xml) — columns A:P
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "27ed5c2c-acc1-4e8a-a34b-ed92d1d00444.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "False"
$ws.Range("G2").Value = "27ed5c2c-acc1-4e8a-a34b-ed92d1d00444.3f9bb167b8acd64d748426ebd885572042073a9a.de-de.xlf"
$ws.Range("H2").Value = "2016-08-28 04:39:28"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I2").Value = "'"
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = "'"
$ws.Range("J2").Style = "Normal"
$ws.Range("K2").Value = "0001-01-01 00:00:00"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L2").Value = "'"
$ws.Range("L2").Style = "Normal"
$ws.Range("M2").Value = "True"
$ws.Range("N2").Value = "'"
$ws.Range("N2").Style = "Normal"
$ws.Range("O2").Value = "False"
$ws.Range("P2").Value = "'"
$ws.Range("P2").Style = "Normal"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f74ea00f5892552ff51f6e0afefea33abf807e6/e2e/27ed5c2c-acc1-4e8a-a34b-ed92d1d00444.md", "", "", "27ed5c2c-acc1-4e8a-a34b-ed92d1d00444.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f74ea00f5892552ff51f6e0afefea33abf807e6/e2e/93ca8fcb-a7e6-49dd-96f9-8a8ef9039be9.md", "", "", "93ca8fcb-a7e6-49dd-96f9-8a8ef9039be9.md")

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))

Write-Host "Generate Report for Handoff - done"
